$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The affected cells store plain numeric-looking values as TEXT
# (shared strings), matching how the source workbook was authored.
# Force text formatting first so Excel doesn't auto-convert the
# typed values into numbers.
$cells = @("B11","C11","D11","B12","C12","D12","B14","C14","D14")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Enterprises density (per 1000 people) - row 11
$ws.Range("B11").Value = "12.91"
$ws.Range("C11").Value = "9.53"
$ws.Range("D11").Value = "22.44"

# Employment (% of total) - row 12
$ws.Range("B12").Value = "10.48"
$ws.Range("C12").Value = "33.39"
$ws.Range("D12").Value = "43.87"

# Enterprises (% of total) - row 14
$ws.Range("B14").Value = "28.76"
$ws.Range("C14").Value = "37.94"
$ws.Range("D14").Value = "89.32"
